# Auto-upload VRF Excel file
# Add a new "cat" worksheet at the end of the workbook with the standard
# VRF model header row (Outdoor/Indoor model, quantity, serial columns)
# used by the other sheets in this workbook.

$wb = $excel.ActiveWorkbook

# Adding/activating a sheet changes the active sheet - remember the
# current one so it can be restored at the end.
$origActiveSheetName = $wb.ActiveSheet.Name

# Copy the last existing sheet so the new sheet naturally inherits this
# workbook's sheet setup (outline props, page margins, header style),
# then place the copy right after it (i.e. at the very end).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "cat"

# Drop everything but the header row that was copied over.
$newSheet.Rows("2:1000").Delete()

# (Re)write the header row text.
$newSheet.Range("A1").Value = "Outdoor Model"
$newSheet.Range("B1").Value = "Outdoor Quantity"
$newSheet.Range("C1").Value = "Outdoor Serial(s)"
$newSheet.Range("D1").Value = "Indoor Model"
$newSheet.Range("E1").Value = "Indoor Quantity"
$newSheet.Range("F1").Value = "Indoor Serial(s)"

# Reset selection to A1 (default for a new sheet) instead of the
# inherited selection from the copied sheet.
[void]$newSheet.Range("A1").Select()

# Restore the original active sheet so this edit only adds "cat".
$wb.Worksheets.Item($origActiveSheetName).Activate()
